$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create new row 53 by copying row 52 (data + formatting) and inserting it
# immediately below, pushing nothing (row 53 was previously blank/out of range).
$ws.Rows(52).Copy()
$ws.Rows(53).Insert(-4121)

# Set the label for the new run
$ws.Cells.Item(53, 2).Value2 = "Baseline 2010-18_C156"

# Fill in the regression-test values for the new row
$ws.Cells.Item(53, 4).Value2 = 1208.5438095555555
$ws.Cells.Item(53, 5).Value2 = 1901.5157334444443
$ws.Cells.Item(53, 6).Value2 = 0.97970299999999988
$ws.Cells.Item(53, 7).Value2 = 280.33542888888883
$ws.Cells.Item(53, 8).Value2 = 9.775355222222224
$ws.Cells.Item(53, 9).Value2 = 5.3172314444444444
$ws.Cells.Item(53, 10).Value2 = 8.145128999999999
$ws.Cells.Item(53, 11).Value2 = 645.97183577777787
$ws.Cells.Item(53, 12).Value2 = 83.47062044444445
$ws.Cells.Item(53, 13).Value2 = 1460.2614338888889
$ws.Cells.Item(53, 14).Value2 = 1208.0519340000001
$ws.Cells.Item(53, 15).Value2 = 4662.6060926666669
$ws.Cells.Item(53, 16).Value2 = 27227.338324888889
$ws.Cells.Item(53, 17).Value2 = -0.5663084444444445
$ws.Cells.Item(53, 18).Value2 = -0.00015933333333333332

# Highlight the three "new" summary columns (D, M, N) in yellow, matching
# the convention used for the most-recent run in this sheet.
$ws.Cells.Item(53, 4).Interior.Color = 65535
$ws.Cells.Item(53, 13).Interior.Color = 65535
$ws.Cells.Item(53, 14).Interior.Color = 65535

# Update the active selection to reflect where the user left off editing.
[void]$ws.Range("S54").Select()
